# Applies two changes described by the commit diff:
#   1. Collapse the three detailed "CORE COMPETENCIES" bullet paragraphs into
#      a single summary line of just the three category names.
#   2. Append a new "TECHNICAL SKILLS" section (Heading 2 + three detail
#      lines) after the "Led multi-million dollar..." bullet, mirroring the
#      abbreviated category/detail pairs that used to live under
#      CORE COMPETENCIES.

$d = $word.ActiveDocument
$bullet = [char]0x2022

# === Change 1: CORE COMPETENCIES -> single condensed paragraph =============
# Paragraph 6 holds the first (Statistical Analysis & Machine Learning...)
# bullet; paragraphs 7 and 8 hold the Big Data and Data Visualization
# bullets respectively. Rewrite paragraph 6 in place, then remove the two
# now-redundant paragraphs that follow it.
$d.Paragraphs(6).Range.Text = "Statistical Analysis & Machine Learning $bullet Big Data & Data Engineering $bullet Data Visualization & Reporting"
$d.Paragraphs(7).Range.Delete() | Out-Null
$d.Paragraphs(7).Range.Delete() | Out-Null

# === Change 2: insert TECHNICAL SKILLS section ==============================
# Locate the "Led multi-million dollar..." bullet (last bullet of the
# "Statistical Analysis & Research" subsection) and insert the new section
# right after it, before the closing "For a more detailed..." paragraph.
$oldText = "Led multi-million dollar research projects involving sensitive consumer data with privacy compliance"

$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*$oldText*") {
        $targetIdx = $i
        break
    }
}

$newText = $oldText + "^pTECHNICAL SKILLS^pSTATISTICAL ANALYSIS & MACHINE LEARNING Advanced Statistical Modeling; Predictive Analytics; Data Mining; Machine Learning^pBIG DATA & DATA ENGINEERING Big Data Processing; Data Warehousing; Cloud Platforms; Data Pipeline Optimization^pDATA VISUALIZATION & REPORTING Data Visualization; Geospatial Analysis; Interactive Dashboards; Business Intelligence"

$d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# The first of the newly-inserted paragraphs ("TECHNICAL SKILLS") becomes a
# Heading 2, matching the other section headers (CORE COMPETENCIES,
# PROFESSIONAL EXPERIENCE, KEY ACHIEVEMENTS AND IMPACT, ...).
$d.Paragraphs($targetIdx + 1).Style = "Heading 2"
